$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 222, shifting existing rows 222:240 down to 223:241.
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new record's data.
$ws.Cells.Item(222, 1).Value = 10
$ws.Cells.Item(222, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(222, 3).Value = "La Araucanía"
$ws.Cells.Item(222, 4).Value = 44746
$ws.Cells.Item(222, 5).Value = 9
$ws.Cells.Item(222, 6).Value = 100112052
$ws.Cells.Item(222, 7).Value = "Albahaca"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 55
$ws.Cells.Item(222, 11).Value = 6000
$ws.Cells.Item(222, 12).Value = 6000
$ws.Cells.Item(222, 13).Value = 6000
$ws.Cells.Item(222, 14).Value = "$/paquete"
$ws.Cells.Item(222, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(222, 16).Value = 6000
$ws.Cells.Item(222, 17).Value = 1
$ws.Cells.Item(222, 18).Value = "Hortaliza"
